# Apply updated symbol-list values (Fri Jan 27 15:03:33 UTC 2023 run).
# Cells are stored as text in the workbook, so numeric-looking values are
# written with a leading quote to force a Text cell (matches original inlineStr
# typing) instead of letting Excel auto-convert them to Number/Percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'304.03"
$ws.Range("E2").Value = "'-0.92%"
$ws.Range("G2").Value = "'15"

# Row 3
$ws.Range("D3").Value = "'35.75"
$ws.Range("E3").Value = "'-0.47%"
$ws.Range("G3").Value = "'15"

# Row 4
$ws.Range("D4").Value = "'5.028"
$ws.Range("E4").Value = "'-0.98%"
$ws.Range("G4").Value = "'15"

# Row 5
$ws.Range("D5").Value = "'0.08015"
$ws.Range("E5").Value = "'-0.85%"
$ws.Range("G5").Value = "'15"

# Row 6
$ws.Range("D6").Value = "'1.852"
$ws.Range("E6").Value = "'-4.32%"
$ws.Range("G6").Value = "'15"

# Row 7
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'7.760"
$ws.Range("E7").Value = "'-0.93%"
$ws.Range("G7").Value = "'15"

# Row 8
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9232"
$ws.Range("E8").Value = "'-1.78%"
$ws.Range("G8").Value = "'15"

# Row 9
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1277"
$ws.Range("E9").Value = "'-4.02%"
$ws.Range("G9").Value = "'15"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1887"
$ws.Range("E10").Value = "'-1.16%"
$ws.Range("G10").Value = "'15"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09001"
$ws.Range("E11").Value = "'-2.88%"
$ws.Range("G11").Value = "'15"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03442"
$ws.Range("E12").Value = "'-2.29%"
$ws.Range("G12").Value = "'15"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09874"
$ws.Range("E13").Value = "'-0.09%"
$ws.Range("G13").Value = "'15"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001402"
$ws.Range("E14").Value = "'-2.02%"
$ws.Range("G14").Value = "'15"

# Row 15
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006237"
$ws.Range("E15").Value = "'6.10%"
$ws.Range("G15").Value = "'15"

# Row 16
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.863"
$ws.Range("E16").Value = "'7.17%"
$ws.Range("G16").Value = "'15"

# Row 17
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.115"
$ws.Range("E17").Value = "'-1.30%"
$ws.Range("G17").Value = "'15"

# Row 18
$ws.Range("D18").Value = "'3.313"
$ws.Range("E18").Value = "'10.47%"
$ws.Range("G18").Value = "'15"

# Row 19
$ws.Range("D19").Value = "'0.3406"
$ws.Range("E19").Value = "'-1.42%"
$ws.Range("G19").Value = "'15"

# Row 20
$ws.Range("D20").Value = "'0.1339"
$ws.Range("E20").Value = "'-0.53%"
$ws.Range("G20").Value = "'15"

# Row 21
$ws.Range("D21").Value = "'4.796"
$ws.Range("E21").Value = "'-7.83%"
$ws.Range("G21").Value = "'15"

# Row 22
$ws.Range("D22").Value = "'0.2336"
$ws.Range("E22").Value = "'-11.06%"
$ws.Range("G22").Value = "'15"

# Row 23
$ws.Range("D23").Value = "'0.04372"
$ws.Range("E23").Value = "'-0.77%"
$ws.Range("G23").Value = "'15"

# Row 24
$ws.Range("D24").Value = "'0.001231"
$ws.Range("E24").Value = "'-0.83%"
$ws.Range("G24").Value = "'15"

# Row 25
$ws.Range("E25").Value = "'2.44%"
$ws.Range("G25").Value = "'15"

# Row 26
$ws.Range("G26").Value = "'15"

# Row 27
$ws.Range("D27").Value = "'0.0001302"
$ws.Range("E27").Value = "'-0.27%"
$ws.Range("G27").Value = "'15"

# Row 28
$ws.Range("E28").Value = "'41.62%"
$ws.Range("G28").Value = "'15"

# Row 29
$ws.Range("G29").Value = "'15"

# Row 30
$ws.Range("G30").Value = "'15"

# Row 31
$ws.Range("G31").Value = "'15"

# Row 32
$ws.Range("G32").Value = "'15"

# Row 33
$ws.Range("G33").Value = "'15"

# Row 34
$ws.Range("G34").Value = "'15"

# Row 35
$ws.Range("G35").Value = "'15"

# Row 36
$ws.Range("G36").Value = "'15"

# Row 37
$ws.Range("G37").Value = "'15"

# Row 38
$ws.Range("G38").Value = "'15"

# Row 39
$ws.Range("D39").Value = "'0.01942"
$ws.Range("E39").Value = "'-2.41%"
$ws.Range("G39").Value = "'15"

# Row 40
$ws.Range("D40").Value = "'0.05115"
$ws.Range("E40").Value = "'2.16%"
$ws.Range("G40").Value = "'15"

# Row 41
$ws.Range("D41").Value = "'0.007565"
$ws.Range("E41").Value = "'-0.77%"
$ws.Range("G41").Value = "'15"

# Row 42
$ws.Range("D42").Value = "'0.01009"
$ws.Range("E42").Value = "'-10.54%"
$ws.Range("G42").Value = "'15"

# Row 43
$ws.Range("D43").Value = "'0.1353"
$ws.Range("E43").Value = "'-1.76%"
$ws.Range("G43").Value = "'15"

# Row 44
$ws.Range("D44").Value = "'0.002113"
$ws.Range("E44").Value = "'0.20%"
$ws.Range("G44").Value = "'15"

# Row 45
$ws.Range("D45").Value = "'0.009855"
$ws.Range("E45").Value = "'-13.27%"
$ws.Range("G45").Value = "'15"

# Row 46
$ws.Range("D46").Value = "'0.00006172"
$ws.Range("E46").Value = "'-3.98%"
$ws.Range("G46").Value = "'15"

# Row 47
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.31%"
$ws.Range("G47").Value = "'15"

# Row 48
$ws.Range("D48").Value = "'63.69"
$ws.Range("E48").Value = "'0.20%"
$ws.Range("G48").Value = "'15"

# Row 49
$ws.Range("D49").Value = "'0.001249"
$ws.Range("E49").Value = "'4.49%"
$ws.Range("G49").Value = "'15"

# Row 50
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.31%"
$ws.Range("G50").Value = "'15"

# Row 51
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.31%"
$ws.Range("G51").Value = "'15"
